$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

# Text columns (Date, Time, Weekday, Week) - force text type so values like
# "2023-06-08" and "23" aren't auto-converted to dates/numbers, then clear
# the temporary formatting so no extra style is left on the cell.
$textCols = 1,2,3,4
foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "2023-06-08"
$ws.Cells.Item($row, 2).Value = "12:51:01"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "23"

foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).ClearFormats()
}

# Numeric columns (city resale counts)
$ws.Cells.Item($row, 5).Value = 118796
$ws.Cells.Item($row, 6).Value = 134371
$ws.Cells.Item($row, 7).Value = 160178
$ws.Cells.Item($row, 8).Value = 131151
$ws.Cells.Item($row, 9).Value = 175543
$ws.Cells.Item($row, 10).Value = 113018
$ws.Cells.Item($row, 11).Value = 200990
$ws.Cells.Item($row, 12).Value = 220957
$ws.Cells.Item($row, 13).Value = 172774
$ws.Cells.Item($row, 14).Value = 120047
$ws.Cells.Item($row, 15).Value = 38594
$ws.Cells.Item($row, 16).Value = 34475
$ws.Cells.Item($row, 17).Value = 50802
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36610
$ws.Cells.Item($row, 20).Value = -1
